$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3:M6").EntireRow.Delete()
$ws.Range("A2").Value = 'Name  Pos  Team  Age  DOB  POB  Nationality  Bats  Throws  Height  Weight  Salary  Babb, Tristan  SP  Green Bay Packers  28  09/04/1973  Phoenix  USA  Right  Right  6'' 2" 210 lbs  -  Baez, Alberto  RP  Washington Commanders  28  10/07/1973  Vacaville  USA  Right  Right  6'' 6" 225 lbs  $250,000  Bailey, Jim  RP  Tampa Bay Buccaneers  26  08/11/1975  Bristow  USA  Left  Left  6'' 0" 195 lbs  -  Bailey, Randy  LF  Tampa Bay Buccaneers  39  07/21/1962  Mesquite  USA  Left  Left  6'' 2" 245 lbs  $244,000  Baker, Chris  3B  Philadelphia Eagles  29  08/29/1972  Melrose  USA  Right  Right  6'' 2" 210 lbs  $7,000,000  Baker, Dustin  3B  New England Patriots  25  06/03/1976  New York  USA  Right  Right  6'' 2" 200 lbs  $1,600,000  Baker, Jason  RF  Detroit Lions  19  09/26/1982  North Creek  USA  Left  Left  6'' 1" 200 lbs  -  Baker, Matt  SS  Tennessee Titans  24  06/26/1977  West Covina  USA  Right  Right  6'' 0" 195 lbs  $244,000  Baker, Shane  CF  Carolina Panthers  39  03/31/1963  Severna Park  USA  Right  Right  6'' 2" 230 lbs  $5,500,000  Baker, Tommy  1B  San Diego Chargers  26  09/25/1975  Waimanalo  USA  Left  Left  6'' 7" 295 lbs  -  Balazs, Jason  LF  Houston Texans  29  10/31/1972  North Shore  USA  Right  Left  5'' 10"  220 lbs  $470,000  Baldwin, Jadon  RP  San Diego Chargers  30  10/31/1971  Spring  USA  Right  Right  6'' 6" 235 lbs  -  Barahona, Rafael  RP  Tampa Bay Buccaneers  26  11/23/1975  Naguanagua  VEN  Right  Right  6'' 8" 240 lbs  -  Barajas, Alvaro  RP  Houston Texans  27  05/25/1974  Santiago  DOM  Right  Right  5'' 9" 180 lbs  -  Barajas, Victor  3B  Washington Commanders  25  04/17/1976  Santo Domingo  DOM  Right  Right  6'' 4" 205 lbs  $1,600,000  Barbosa, Jesus  2B  Carolina Panthers  16  09/21/1985  San Francisco de Macorís  DOM  Right  Right  5'' 10"  175 lbs  -  Barcenas, Bill  SP  Cleveland Browns  16  10/11/1985  Tlalnepantla  MEX  Right  Right  6'' 4" 190 lbs  -  Barclay, Ryan  RP  Houston Texans  25  03/02/1977  Douglas  USA  Left  Left  6'' 2" 210 lbs  $380,000  Barendse, Adan  RP  Washington Commanders  24  05/31/1977  South Windsor  USA  Switch  Right  6'' 0" 195 lbs  $800,000  Barnes, Keith  CL  Houston Texans  31  09/24/1970  Altoona  USA  Right  Right  6'' 0" 205 lbs  $4,000,000  Barrera, Jose  RP  Green Bay Packers  29  09/01/1972  Valencia  VEN  Right  Right  6'' 0" 205 lbs  $244,000  Barrera, Wilson  C  Detroit Lions  27  07/26/1974  Salcedo  DOM  Right  Right  5'' 11"  200 lbs  $2,800,000  Barrett, Matt  RF  Seattle Seahawks  33  05/03/1968  Yonkers  USA  Right  Left  6'' 4" 255 lbs  -  Barriga, Alex  CL  Kansas City Chiefs  23  04/11/1979  Wading River  USA  Right  Right  5'' 11"  190 lbs  $244,000  Barrosa, Pedro  1B  Chicago Bears  27  09/29/1974  Valencia  VEN  Right  Right  6'' 5" 245 lbs  $800,000  Barroso, Pedro  LF  Washington Commanders  28  09/27/1973  Santo Domingo  DOM  Right  Right  6'' 1" 235 lbs  -  Bartlett, Josh  SS  Seattle Seahawks  37  12/23/1964  Lake Forest  USA  Right  Right  6'' 1" 205 lbs  $590,000  Bartmess, Nick  RP  Detroit Lions  33  01/28/1969  Lowell  USA  Left  Left  5'' 10"  195 lbs  -  Bartnick, Pat  2B  St. Louis Rams  31  08/11/1970  Gun Barrel City  USA  Right  Right  6'' 4" 230 lbs  $620,000  Bass, Isaiah  SP  Carolina Panthers  27  03/04/1975  Dallas  USA  Right  Right  6'' 1" 205 lbs  $306,000  Bass, Khalil  1B  San Francisco 49ers  24  09/03/1977  Parkway  USA  Left  Left  6'' 1" 210 lbs  $2,000,000  Battllori, Renault  SP  New York Giants  36  02/13/1966  Venice  ITA  Left  Left  6'' 2" 215 lbs  -  Baumann, Tommy  3B  Detroit Lions  17  07/10/1984  Douglasville  USA  Switch  Right  6'' 1" 180 lbs  $244,000  Baumgartner, Nick  SS  Atlanta Falcons  29  08/28/1972  Seattle  USA  Right  Right  6'' 2" 200 lbs  $250,000  Bayne, Brian  RP  Oakland Raiders  29  04/19/1972  Caribou  USA  Right  Left  6'' 0" 200 lbs  -  Bazaldua, Tony  RP  Miami Dolphins  25  01/21/1977  Maracaibo  VEN  Left  Left  6'' 1" 200 lbs  $460,000  Beardslee, John  CF  New York Giants  18  07/17/1983  Baltimore  USA  Right  Right  5'' 11"  170 lbs  $244,000  Becerra, Omar  1B  Baltimore Ravens  29  01/21/1973  Bonao  DOM  Left  Right  6'' 0" 235 lbs  $650,000  Beck, Tommy  1B  Washington Commanders  28  03/29/1974  North Mankato  USA  Right  Right  6'' 1" 220 lbs  $244,000  Bedney, Jim  1B  Cincinnati Bengals  30  02/05/1972  Pottsville  USA  Left  Left  6'' 6" 255 lbs  $8,300,000  Beebe, John  1B  Oakland Raiders  27  08/02/1974  Pomona  USA  Right  Right  6'' 2" 240 lbs  -  Beery, Mike  LF  Denver Broncos  37  01/11/1965  Van Buren  USA  Left  Left  6'' 4" 245 lbs  $276,000  Bejarano, Jose  RP  Chicago Bears  27  06/07/1974  Levittown  PUR  Right  Right  6'' 0" 200 lbs  -  Belford, Landon  SS  Miami Dolphins  26  03/13/1976  Baltimore  USA  Right  Right  6'' 0" 210 lbs  -  Bell, Amari  RP  Pittsburgh Steelers  28  05/21/1973  Nashua  USA  Right  Right  6'' 2" 210 lbs  $434,000  Belmore, Del  LF  Oakland Raiders  28  02/21/1974  Hamilton  CAN  Right  Right  6'' 3" 205 lbs  $5,700,000  Beltran, Marcos  RP  Cincinnati Bengals  28  03/31/1974  Santo Domingo  DOM  Right  Right  6'' 2" 210 lbs  $250,000  Bencosme, Luis  LF  San Francisco 49ers  27  09/14/1974  San Cristóbal  DOM  Left  Left  6'' 4" 250 lbs  -  Benedict, Nick  RF  Philadelphia Eagles  24  07/16/1977  Magna  USA  Right  Right  6'' 1" 210 lbs  $1,900,000  Benitez, Humberto  LF  St. Louis Rams  26  01/12/1976  Santa Teresa  VEN  Right  Right  6'' 2" 235 lbs  -  Benitez, Luis  RF  Detroit Lions  36  06/15/1965  Santo Domingo  DOM  Right  Right  6'' 6" 250 lbs  $1,840,000  Benjamine, Josh  3B  Jacksonville Jaguars  26  12/23/1975  Sherwood Manor  USA  Right  Right  6'' 4" 215 lbs  $9,800,000  Bennett, Elijah  CF  Washington Commanders  27  02/26/1975  Mirabel  CAN  Left  Right  6'' 1" 195 lbs  $284,000  Benning, Paul  1B  Dallas Cowboys  31  08/26/1970  Vancouver  USA  Left  Left  6'' 4" 265 lbs  -  Benson, Ryan  SP  New York Giants  35  03/17/1967  Calexico  USA  Switch  Right  6'' 5" 230 lbs  $3,200,000  Benson, Steve  CF  Dallas Cowboys  32  01/01/1970  Mountain View  USA  Right  Right  6'' 1" 215 lbs  $244,000  Benvenuto, C.J.  RP  Miami Dolphins  23  12/06/1978  Elk Grove  USA  Right  Right  6'' 0" 195 lbs  $1,000,000  Berg, Chris  C  New York Jets  25  06/13/1976  Jacksonville  USA  Left  Right  6'' 6" 245 lbs  $1,700,000  Bernacki, Paul  RF  Kansas City Chiefs  28  04/12/1974  Sand Springs  USA  Left  Left  6'' 0" 195 lbs  $2,880,000  Bernal, Ricky  3B  New England Patriots  27  01/22/1975  Germantown  USA  Right  Right  6'' 0" 185 lbs  $244,000  Betancur, Johnny  1B  Cleveland Browns  32  03/21/1970  Santo Domingo  DOM  Right  Right  6'' 4" 230 lbs  $9,100,000  Bettiza, Norbert  C  San Diego Chargers  16  09/10/1985  Portogruaro  ITA  Right  Right  6'' 1" 195 lbs  -  Bice, Justin  SP  New Orleans Saints  19  03/02/1983  Opportunity  USA  Right  Right  6'' 4" 190 lbs  -  Biggs, Jeff  RP  Arizona Cardinals  27  08/13/1974  Falmouth  USA  Left  Left  6'' 1" 195 lbs  -  Bisceglia, Dave  SS  Carolina Panthers  34  06/08/1967  New York  USA  Right  Right  6'' 0" 200 lbs  $860,000  Bishop, Jonathan  1B  San Diego Chargers  21  07/01/1980  Hamilton  CAN  Switch  Right  6'' 2" 235 lbs  -  Black, Dylan  RF  Chicago Bears  33  08/27/1968  Woodland Heights  USA  Left  Left  6'' 3" 225 lbs  $2,800,000  Blackburn, Steve  C  Chicago Bears  33  12/04/1968  Costa Mesa  USA  Right  Right  6'' 0" 215 lbs  $244,000  Blake, Sean  2B  Philadelphia Eagles  24  10/09/1977  Ross Township  USA  Left  Right  6'' 1" 210 lbs  $800,000  Blakney, Mike  CF  Denver Broncos  30  05/08/1971  Los Angeles  USA  Left  Left  6'' 1" 200 lbs  $260,000  Blanchette, Warren  RP  Indianapolis Colts  24  12/22/1977  Hampton  USA  Switch  Right  6'' 1" 200 lbs  -  Bland, Brian  RP  Miami Dolphins  28  06/05/1973  Sacramento  USA  Right  Right  6'' 1" 200 lbs  -  Blank, Gabe  SS  Cleveland Browns  27  01/29/1975  Lochbuie  USA  Right  Right  6'' 1" 195 lbs  $580,000  Blankenship, Andy  RP  Detroit Lions  21  10/01/1980  Lancaster  USA  Switch  Left  6'' 5" 215 lbs  -  Bojorquez, Sergio  RP  Tampa Bay Buccaneers  32  11/19/1969  San Cristóbal  VEN  Switch  Right  6'' 0" 200 lbs  -  Bolen, Mike  RP  Dallas Cowboys  35  05/28/1966  Spencer  USA  Right  Right  6'' 1" 210 lbs  -  Bonaddio, Zach  RP  Tampa Bay Buccaneers  29  01/07/1973  Sydney  AUS  Right  Left  6'' 2" 205 lbs  -  Bonilla, Carlos  LF  Detroit Lions  21  05/16/1980  Río Verde Arriba  DOM  Left  Left  5'' 10"  200 lbs  -  Bonnar, Edison  RP  New Orleans Saints  31  09/13/1970  Elmira  CAN  Right  Right  5'' 10"  185 lbs  -  Bordelon, Eric  RP  Jacksonville Jaguars  28  05/03/1973  Easley  USA  Right  Right  6'' 8" 245 lbs  $402,000  Borho, Ian  LF  New Orleans Saints  28  10/03/1973  Olympia  USA  Left  Left  6'' 0" 220 lbs  -  Borjas, Cesar  1B  Carolina Panthers  22  07/31/1979  Canton  USA  Right  Right  6'' 2" 220 lbs  $244,000  Boros, Dave  SS  Carolina Panthers  19  04/26/1982  Logan  USA  Right  Right  6'' 0" 175 lbs  -  Botello, Marco  RP  Oakland Raiders  28  04/12/1974  Barquisimeto  VEN  Right  Right  6'' 6" 225 lbs  $420,000  Boubret, Fabrício  RP  Minnesota Vikings  32  08/01/1969  Maracaibo  VEN  Right  Right  5'' 10"  195 lbs  $420,000  Bourke, Dan  RP  San Diego Chargers  24  10/27/1977  Lake City  USA  Left  Left  6'' 4" 220 lbs  -  Bowen, Bobby  LF  Green Bay Packers  35  05/10/1966  Ramapo  USA  Left  Left  6'' 4" 260 lbs  -  Bowie, Nate  RF  Jacksonville Jaguars  22  10/24/1979  Peabody  USA  Left  Left  6'' 5" 225 lbs  $244,000  Boyer, Barry  RP  Detroit Lions  28  09/18/1973  Ortonville  USA  Right  Right  6'' 0" 200 lbs  $3,240,000  Brackman, Bob  2B  Houston Texans  24  10/08/1977  Rancho Cordova  USA  Right  Right  6'' 0" 190 lbs  $1,800,000  Brady, David  CF  Arizona Cardinals  24  11/10/1977  Atlantic Beach  USA  Right  Right  6'' 2" 220 lbs  -  Bragg, Chris  RP  Baltimore Ravens  26  04/14/1976  Garnet Valley  USA  Left  Right  6'' 2" 205 lbs  -  Braithwaite, Joel  RP  Cincinnati Bengals  33  09/02/1968  Hurst  USA  Left  Left  6'' 0" 205 lbs  $890,000  Bresler, Aaron  RP  Philadelphia Eagles  23  02/28/1979  Los Angeles  USA  Left  Left  6'' 2" 200 lbs  -  Bretz, Tyler  RP  Seattle Seahawks  28  12/14/1973  Temecula  USA  Right  Right  5'' 10"  190 lbs  -  Briones, Francisco  SS  Philadelphia Eagles  22  10/15/1979  Caracas  VEN  Right  Right  6'' 4" 210 lbs  -  Brito, Manny  C  San Diego Chargers  24  12/15/1977  Santiago  DOM  Left  Right  6'' 3" 240 lbs  $410,000  Britt, Jody  LF  Denver Broncos  24  04/16/1977  Endicott  USA  Right  Right  6'' 6" 235 lbs  $900,000  Brody, Danny  RP  Minnesota Vikings  35  01/29/1967  Lakeview  USA  Left  Left  6'' 0" 210 lbs  -  Brogle, Mark  LF  Green Bay Packers  32  01/17/1970  Ketchikan  USA  Right  Left  6'' 1" 235 lbs  $354,000  Brooks, Josh  RP  Miami Dolphins  35  10/30/1966  Mineola  USA  Right  Left  6'' 2" 215 lbs  -  Broulik, Sam  RP  Free Agent  17  11/20/1984  West Covina  USA  Left  Left  6'' 3" 190 lbs  -  Brown, Alex  RF  Seattle Seahawks  27  10/15/1974  New York  USA  Left  Right  6'' 3" 210 lbs  $8,300,000  Brown, Jerry  3B  Cleveland Browns  30  06/13/1971  Laurens  USA  Right  Right  6'' 4" 235 lbs  $434,000  Brown, Joe  LF  San Francisco 49ers  26  12/12/1975  Waterville  USA  Right  Left  6'' 2" 220 lbs  $440,000  Brown, Josh  RP  Carolina Panthers  24  05/15/1977  Gilbert  USA  Left  Left  6'' 0" 200 lbs  $750,000  Brown, Mike  RP  Tampa Bay Buccaneers  31  07/31/1970  Chicago  USA  Right  Right  6'' 4" 220 lbs  -  Browning, Josh  RP  Washington Commanders  31  04/12/1971  Irving  USA  Left  Right  6'' 5" 230 lbs  -  Brunetti, Darrell  RP  Buffalo Bills  33  06/22/1968  Watsonville  USA  Right  Right  6'' 4" 220 lbs  -  Bryant, Mike  RP  Baltimore Ravens  24  10/14/1977  Seven Corners  USA  Left  Left  6'' 4" 220 lbs  $350,000  Buchanan, Jordan  RP  Indianapolis Colts  28  01/06/1974  El Paso  USA  Right  Left  6'' 2" 210 lbs  -  Budd, Justin  RP  Free Agent  20  04/03/1982  Jackson  USA  Right  Right  6'' 2" 185 lbs  -  Buesing, Pete  RP  Indianapolis Colts  19  02/12/1983  Orange Cove  USA  Switch  Right  6'' 2" 180 lbs  -  Burke, Bobby  RP  Miami Dolphins  25  03/29/1977  Hutchison  USA  Left  Left  6'' 1" 200 lbs  -  Burns, Dale  RP  St. Louis Rams  29  12/27/1972  Greenville  USA  Right  Right  6'' 1" 205 lbs  -  Burrows, Andy  RF  Cincinnati Bengals  30  09/02/1971  Cortlandt  USA  Right  Right  6'' 1" 215 lbs  $2,080,000  Burruel, Juan  SP  Arizona Cardinals  16  01/31/1986  Bogota  COL  Left  Left  6'' 0" 165 lbs  -  Burtrum, Heath  CF  Atlanta Falcons  27  12/27/1974  Brown City  USA  Right  Right  6'' 2" 205 lbs  $600,000  Buss, Jesse  SS  San Francisco 49ers  27  09/21/1974  Louisville  USA  Right  Right  6'' 0" 175 lbs  $11,500,000  Bustamante, Luis  SP  Atlanta Falcons  32  10/25/1969  San Juan  PUR  Right  Right  6'' 4" 225 lbs  $650,000  Butterell, Warren  RP  New York Giants  35  03/31/1967  Adelaide  AUS  Right  Right  6'' 1" 205 lbs  $244,000  Bynum, Nate  LF  Chicago Bears  28  07/31/1973  San Diego  USA  Right  Left  6'' 4" 230 lbs  $2,600,000'
